$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Raw data sheets (NaCl, CaCl2, SiO2):
#    - the second "Concentrate" sample per time point is relabelled
#      "Retentate"
#    - the last time point (previously 4.5) is corrected to 4
#    - leave each sheet's own selection parked on D15
# ---------------------------------------------------------------------------
$rawSheetNames = @("NaCl", "CaCl2", "SiO2")
foreach ($name in $rawSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("D4").Value = "Retentate"
    $ws.Range("D7").Value = "Retentate"
    $ws.Range("D10").Value = "Retentate"
    $ws.Range("D13").Value = "Retentate"
    $ws.Range("D15").Value = "Retentate"

    $ws.Range("A17").Value = 4

    [void]$ws.Range("D15").Select()
}

# ---------------------------------------------------------------------------
# 2) "_check" sheets (NaCl_Check, CaCl2_check, SiO2_check):
#    same time-point correction (row 19 there, since these sheets carry
#    three header rows before the data starts), plus their own parked
#    selections.
# ---------------------------------------------------------------------------
$checkSelections = @{
    "NaCl_Check"  = "F5"
    "CaCl2_check" = "A20"
    "SiO2_check"  = "B22"
}
foreach ($name in $checkSelections.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A19").Value = 4
    [void]$ws.Range($checkSelections[$name]).Select()
}

# ---------------------------------------------------------------------------
# 3) NaCl ends up the active sheet/tab, selection resting on D15.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NaCl")
$ws.Activate()
[void]$ws.Range("D15").Select()
